$d = $word.ActiveDocument

# Locate the start of the sentence that contains the renewal deadline text.
$anchor = $d.Content
$anchor.Find.Execute("les équipes seront constituées", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base = $anchor.Start

# Drop Word's existing "last edit" bookmark - it currently sits in the middle
# of " la cotisat|ion ..." and would otherwise block those two fragments from
# re-joining into a single run below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Character offsets (relative to $base) of the internal boundaries we want to
# turn into separate runs, computed against the *original* wording:
#   0   "les équ"
#   7   "ipes seront constituées sans garantie de place à votre niveau habituel. Passé"
#   84  " le vendredi "
#   97  "13"            <- becomes "6"
#   99  " septembre"
#   109 ","
#   110 " la cotisation passe à 90 €."
#   138 end of sentence (protects the trailing space run that follows)
# Bookmarks automatically re-anchor as the text around them is edited, so we
# use them (rather than recomputed arithmetic) for every position we still
# need after a prior edit has shifted the surrounding text.
$splits = 7, 84, 97, 99, 109, 110, 138

$i = 0
foreach ($off in $splits) {
    $i = $i + 1
    $pt = $d.Range($base + $off, $base + $off)
    $d.Bookmarks.Add("zzsplit$i", $pt) | Out-Null
}

# Narrow replace of the deadline day, bounded by the zzsplit3/zzsplit4
# bookmarks so the surrounding runs are left untouched.
$digits = $d.Range($d.Bookmarks.Item("zzsplit3").Range.Start, $d.Bookmarks.Item("zzsplit4").Range.Start)
$digits.Text = "6"

# Re-join " la cotisat" + "ion passe à 90 €." (split apart originally by the
# old _GoBack bookmark) into a single run, bounded by the zzsplit6/zzsplit7
# bookmarks so nothing else is touched. A same-text assignment is a no-op in
# this engine, so round-trip through a placeholder to force the rewrite.
$cotisStart = $d.Bookmarks.Item("zzsplit6").Range.Start
$cotisEnd = $d.Bookmarks.Item("zzsplit7").Range.Start
$cotis = $d.Range($cotisStart, $cotisEnd)
$cotisText = $cotis.Text
$cotis.Text = "x"
$cotis2 = $d.Range($cotisStart, $cotisStart + 1)
$cotis2.Text = $cotisText

# Clean up the temporary bookmarks (removing a bookmark does not merge the
# runs its presence had already split apart).
for ($k = 1; $k -le $splits.Length; $k++) {
    $d.Bookmarks.Item("zzsplit$k").Delete()
}

# Word leaves its "last edit" marker (_GoBack) right where the user's cursor
# ended up - immediately after "les équ".
$gbPoint = $d.Range($base + 7, $base + 7)
$d.Bookmarks.Add("_GoBack", $gbPoint) | Out-Null
